$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.624187469482422
$ws.Range("B1").Value = 6.111722946166992
$ws.Range("C1").Value = 8.243190765380859
$ws.Range("D1").Value = 6.976041793823242
$ws.Range("E1").Value = 2.628762245178223
